$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) for rows 2-7, columns B:G
$data = @(
    @(3.230985683306322, 1.667794583268128, 3.900430680208489,  0.496779210170732, 1, 9.295990156953671),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182),
    @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.997429241610044),
    @(3.230985683306322, 10.29869402782916, 0.8054896365839992,  8.660232485948974, 1, 22.99540183366846)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
